$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I3").Value = "Melanoma risk not associated with age at menarche, age at first use of oral contraceptives, age at first birth, parity, meno status, or HRT use, height, weight, BMI"
$ws.Range("F5").Value = "parity, age at first birth"
$ws.Range("E5").Value = "oral contraceptive use and duration"
$ws.Range("B5").Value = "IJE"
